# Update the public-exposure-sites table with the latest case locations:
#  - Row 2 (the "new" entry) is replaced with the newest exposure site
#    (Emerald / Lakeside Paddle Boats).
#  - Row 3 (the "old" entry) is replaced with the site that was previously
#    listed as "new" (Melbourne / Nandos).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> newest entry ("new")
$ws.Range("A2").Value = "Emerald"
$ws.Range("B2").Value = "Lakeside Paddle Boats, Emerald Lake Park"
$ws.Range("C2").Value = "31/12/20 3:30pm - 5:30pm"
$ws.Range("D2").Value = "Case visited venue"
$ws.Range("E2").Value = "new"

# Row 3 -> previous entry, now marked "old"
$ws.Range("A3").Value = "Melbourne"
$ws.Range("B3").Value = "Nandos  27 Elizabeth Street, Melbourne"
$ws.Range("C3").Value = "01/01/2021 2:00am - 2:30am"
$ws.Range("D3").Value = "Case dined at venue"
$ws.Range("E3").Value = "old"

# Column widths were re-measured (best-fit) by Excel after the content change.
$ws.Columns.Item(1).ColumnWidth = 8.5
$ws.Columns.Item(2).ColumnWidth = 33
$ws.Columns.Item(3).ColumnWidth = 23.5
$ws.Columns.Item(4).ColumnWidth = 15.6666666666667

# Selection moved to B2:B3 with B2 active
$ws.Range("B2:B3").Select() | Out-Null
